$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.429.74'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").Value = '1.621.65'
$ws.Range("E3").Value = '  +0.58%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.38%  '

$ws.Range("E11").Value = '  -0.69%  '

$ws.Range("D12").Value = '1.848.76'
$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("D13").Value = '1.613.72'
$ws.Range("E13").Value = '  -0.28%  '

$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.85'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '236.18'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +6.94%  '

$ws.Range("D18").Value = '26.437.53'
$ws.Range("E18").Value = '  +0.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.84'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.08%  '

$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.19'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.13'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.09'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.50%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("E27").Value = '  +0.65%  '

$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.58'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.04%  '

$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("E31").Value = '  -0.07%  '

$ws.Range("D32").Value = '1.522.73'
$ws.Range("E32").Value = '  +5.26%  '

$ws.Range("E33").Value = '  +1.46%  '

$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.52'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.73%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.567'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0167'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.21%  '

$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("E40").Value = '  +0.71%  '

$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("E42").Value = '  +0.82%  '

$ws.Range("D43").Value = '1.760.58'
$ws.Range("E43").Value = '  +0.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.74'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.763'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.915'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.64'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.94%  '

$ws.Range("E48").Value = '  +1.18%  '

$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("E50").Value = '  +0.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.52'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.28%  '

